# Commit: "commit all crm gui"
# Adds a new "product" worksheet (after Sheet2) with brand/product data,
# makes it the active/selected sheet, and sizes its two columns.

$wb = $excel.ActiveWorkbook

# --- Add the new "product" worksheet after the last existing sheet (Sheet2) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "product"

# --- Populate data (A1:B6) ---
$ws.Cells.Item(1, 1).Value = "Brand name"
$ws.Cells.Item(1, 2).Value = "Product name"

$ws.Cells.Item(2, 1).Value = "iphone"
$ws.Cells.Item(2, 2).Value = "Apple iPhone 14 Plus (Blue, 128 GB)"

$ws.Cells.Item(3, 1).Value = "iphone"
$ws.Cells.Item(3, 2).Value = "Apple iPhone 12 (Blue, 64 GB)"

$ws.Cells.Item(4, 1).Value = "iphone"
$ws.Cells.Item(4, 2).Value = "Apple iPhone 15 (Blue, 128 GB)"

$ws.Cells.Item(5, 1).Value = "iphone"
$ws.Cells.Item(5, 2).Value = "Apple iPhone 14 Plus (Midnight, 128 GB)"

$ws.Cells.Item(6, 1).Value = "iphone"
$ws.Cells.Item(6, 2).Value = "Apple iPhone 15 Plus (Pink, 128 GB)"

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 22
$ws.Columns.Item(2).ColumnWidth = 36

# --- Make "product" the active sheet / selection ---
$ws.Activate() | Out-Null
$ws.Range("J11").Select() | Out-Null

Write-Host "Added product sheet with data"
